$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean sheet - remove the old "random" algorithm data
$ws.Cells.Clear() | Out-Null

# ---- Header row (row 1): bold + centered ----
$ws.Range("A1").Value = "algoritme"
$ws.Range("B1").Value = "average moves 6x6 game #1"
$ws.Range("C1").Value = "amount of measurements"
$ws.Range("E1").Value = "average moves 9x9 game #4"
$ws.Range("F1").Value = "amount of measurements"

$ws.Range("A1:C1").HorizontalAlignment = -4108
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("E1:F1").HorizontalAlignment = -4108
$ws.Range("E1:F1").Font.Bold = $true

# ---- Data rows (rows 2-5): centered ----
$ws.Range("A2").Value = "single step, normal win"
$ws.Range("B2").Value = 597.41656999999998
$ws.Range("C2").Value = "100k"

$ws.Range("A3").Value = "single step, path free win"
$ws.Range("B3").Value = 103.05871999999999
$ws.Range("C3").Value = "1 mln"

$ws.Range("A4").Value = "max step, path free win"
$ws.Range("B4").Value = 101.35699
$ws.Range("C4").Value = "100k"

$ws.Range("A5").Value = "max step, path free win, non recurrent"
$ws.Range("B5").Value = 90.427689999999998
$ws.Range("C5").Value = "100k"

$ws.Range("E5").Value = 6408.165
$ws.Range("F5").Value = "1k"

$ws.Range("A2:C5").HorizontalAlignment = -4108
$ws.Range("E5:F5").HorizontalAlignment = -4108

# ---- Column widths ----
$ws.Range("A1").ColumnWidth = 35.85546875
$ws.Range("B1").ColumnWidth = 26
$ws.Range("C1").ColumnWidth = 24.28515625
$ws.Range("E1").ColumnWidth = 26
$ws.Range("F1").ColumnWidth = 24.28515625

# ---- Page setup ----
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---- Selection ----
$ws.Range("C29").Select() | Out-Null
